$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder header content: move "unitType" column before "phone"/"email" ---
# Before: D=phone, E=email, F=unitType
# After:  D=unitType, E=phone, F=email
$ws.Range("A1:F3").ClearContents()

$ws.Range("A1").Value = "{d.i18n.name}"
$ws.Range("B1").Value = "{d.i18n.address}"
$ws.Range("C1").Value = "{d.i18n.unitName}"
$ws.Range("D1").Value = "{d.i18n.unitType}"
$ws.Range("E1").Value = "{d.i18n.phone}"
$ws.Range("F1").Value = "{d.i18n.email}"

$ws.Range("A2").Value = "{d.contacts[i].name}"
$ws.Range("B2").Value = "{d.contacts[i].address}"
$ws.Range("C2").Value = "{d.contacts[i].unitName}"
$ws.Range("D2").Value = "{d.contacts[I].unitType}"
$ws.Range("E2").Value = "{d.contacts[i].phone}"
$ws.Range("F2").Value = "{d.contacts[i].email}"

$ws.Range("A3").Value = "{d.contacts[i+1].name}"
$ws.Range("B3").Value = "{d.contacts[i+1].address}"
$ws.Range("C3").Value = "{d.contacts[i+1].unitName}"
$ws.Range("D3").Value = "{d.contacts[I+1].unitType}"
$ws.Range("E3").Value = "{d.contacts[i+1].phone}"
$ws.Range("F3").Value = "{d.contacts[i+1].email}"

# --- Column widths: C/D become the old "unitType" width (19.1719), E becomes
#     the old "phone" width (16.5); F is untouched and keeps 19.1719 exactly ---
$ws.Range("C1:D1").ColumnWidth = 18.33
$ws.Range("E1").ColumnWidth = 15.67

# --- Row 4-10 borders: column D now needs the "right edge" border that used
#     to belong to column F, and column E now needs the "left edge" border
#     that used to belong to column A. Copy the formats across so the exact
#     same style indices get reused (no new styles created). ---
$ws.Range("F4:F10").Copy()
$ws.Range("D4:D10").PasteSpecial(-4122)
$ws.Range("A4:A10").Copy()
$ws.Range("E4:E10").PasteSpecial(-4122)
